$wb = $excel.ActiveWorkbook

# ALC row 3
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 44857.145
$ws.Range("J3").Value = 44857.145
$ws.Range("L3").Value = 44857.145
$ws.Range("N3").Value = -45085.145

# ALC row 102
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H102").Value = 44857.145
$ws.Range("J102").Value = 44857.145
$ws.Range("L102").Value = 44857.145
$ws.Range("N102").Value = -51347.145

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 196127.36
$ws.Range("I138").Value = 2225.95
$ws.Range("J138").Value = 303850.38
$ws.Range("K138").Value = 6677.849999999999
$ws.Range("L138").Value = 911551.14
$ws.Range("M138").Value = -1537.849999999999
$ws.Range("N138").Value = -921831.14

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1404.6923
$ws.Range("I2").Value = 1138.5
$ws.Range("J2").Value = 1632.8572
$ws.Range("K2").Value = 1138.5
$ws.Range("L2").Value = 1632.8572
$ws.Range("M2").Value = -1025.5
$ws.Range("N2").Value = -1858.8572

# ARM row 28
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 27754.357
$ws.Range("I28").Value = 4457
$ws.Range("J28").Value = 34108.184
$ws.Range("K28").Value = 4457
$ws.Range("L28").Value = 34108.184
$ws.Range("M28").Value = -4265
$ws.Range("N28").Value = -34492.184

# ARM row 31
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 34069.832
$ws.Range("I31").Value = 2966.6667
$ws.Range("J31").Value = 65173
$ws.Range("K31").Value = 2966.6667
$ws.Range("L31").Value = 65173
$ws.Range("M31").Value = -2672.6667
$ws.Range("N31").Value = -65761

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2234.9656
$ws.Range("I45").Value = 1500
$ws.Range("J45").Value = 3022.4285
$ws.Range("K45").Value = 1500
$ws.Range("L45").Value = 3022.4285
$ws.Range("M45").Value = -1123
$ws.Range("N45").Value = -3776.4285

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1639
$ws.Range("I74").Value = 1466.2222
$ws.Range("J74").Value = 2675.6667
$ws.Range("K74").Value = 1466.2222
$ws.Range("L74").Value = 2675.6667
$ws.Range("M74").Value = -592.2221999999999
$ws.Range("N74").Value = -4423.6667

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1639
$ws.Range("I77").Value = 1466.2222
$ws.Range("J77").Value = 2675.6667
$ws.Range("K77").Value = 7331.111
$ws.Range("L77").Value = 13378.3335
$ws.Range("M77").Value = -2963.111
$ws.Range("N77").Value = -22114.3335

# ARM row 93
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H93").Value = 79224
$ws.Range("J93").Value = 79224
$ws.Range("L93").Value = 79224
$ws.Range("N93").Value = -84216

# ARM row 94
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H94").Value = 30330
$ws.Range("J94").Value = 30330
$ws.Range("L94").Value = 30330
$ws.Range("N94").Value = -32132

# ARM row 99
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 27754.357
$ws.Range("I99").Value = 4457
$ws.Range("J99").Value = 34108.184
$ws.Range("K99").Value = 4457
$ws.Range("L99").Value = 34108.184
$ws.Range("M99").Value = -1462
$ws.Range("N99").Value = -40098.184

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1447.875
$ws.Range("I110").Value = 1226.1428
$ws.Range("K110").Value = 1226.1428
$ws.Range("M110").Value = 818.8571999999999

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1404.6923
$ws.Range("I116").Value = 1138.5
$ws.Range("J116").Value = 1632.8572
$ws.Range("K116").Value = 1138.5
$ws.Range("L116").Value = 1632.8572
$ws.Range("M116").Value = 1155.5
$ws.Range("N116").Value = -6220.8572

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4365.778
$ws.Range("I132").Value = 4374.7856
$ws.Range("J132").Value = 4356.077
$ws.Range("K132").Value = 13124.3568
$ws.Range("L132").Value = 13068.231
$ws.Range("M132").Value = -10594.3568
$ws.Range("N132").Value = -18128.231

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1404.6923
$ws.Range("I3").Value = 1138.5
$ws.Range("J3").Value = 1632.8572
$ws.Range("K3").Value = 1138.5
$ws.Range("L3").Value = 1632.8572
$ws.Range("M3").Value = -1024.5
$ws.Range("N3").Value = -1860.8572

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 960.17645
$ws.Range("I107").Value = 446.375
$ws.Range("K107").Value = 446.375
$ws.Range("M107").Value = 1473.625

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1224.3478
$ws.Range("I58").Value = 789.6667
$ws.Range("K58").Value = 789.6667
$ws.Range("M58").Value = -586.6667

# CRP row 92
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 28601
$ws.Range("J92").Value = 28601
$ws.Range("L92").Value = 28601
$ws.Range("N92").Value = -33593

# CRP row 96
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 57874.668
$ws.Range("J96").Value = 57874.668
$ws.Range("L96").Value = 57874.668
$ws.Range("N96").Value = -63366.668

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1224.3478
$ws.Range("I136").Value = 789.6667
$ws.Range("K136").Value = 2369.0001
$ws.Range("M136").Value = 180.9998999999998

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2509.1667
$ws.Range("I136").Value = 1627.5
$ws.Range("J136").Value = 2950
$ws.Range("K136").Value = 4882.5
$ws.Range("L136").Value = 8850
$ws.Range("M136").Value = 217.5
$ws.Range("N136").Value = -19050

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 4553.3335
$ws.Range("I139").Value = 1372.5
$ws.Range("J139").Value = 30000
$ws.Range("K139").Value = 4117.5
$ws.Range("L139").Value = 90000
$ws.Range("M139").Value = 1022.5
$ws.Range("N139").Value = -100280

# GSM row 39
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 19498.75
$ws.Range("J39").Value = 19498.75
$ws.Range("L39").Value = 19498.75
$ws.Range("N39").Value = -20562.75

# GSM row 95
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 79672
$ws.Range("J95").Value = 79672
$ws.Range("L95").Value = 79672
$ws.Range("N95").Value = -85164

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4653.6553
$ws.Range("I22").Value = 799.7857
$ws.Range("J22").Value = 8250.6
$ws.Range("K22").Value = 799.7857
$ws.Range("L22").Value = 8250.6
$ws.Range("M22").Value = -504.7857
$ws.Range("N22").Value = -8840.6

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 4653.6553
$ws.Range("I27").Value = 799.7857
$ws.Range("J27").Value = 8250.6
$ws.Range("K27").Value = 799.7857
$ws.Range("L27").Value = 8250.6
$ws.Range("M27").Value = -692.7857
$ws.Range("N27").Value = -8464.6

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4617
$ws.Range("I61").Value = 4687.5835
$ws.Range("J61").Value = 4475.8335
$ws.Range("K61").Value = 4687.5835
$ws.Range("L61").Value = 4475.8335
$ws.Range("M61").Value = -4485.5835
$ws.Range("N61").Value = -4879.8335

# LTW row 98
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 58249.5
$ws.Range("J98").Value = 58249.5
$ws.Range("L98").Value = 58249.5
$ws.Range("N98").Value = -64239.5

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4617
$ws.Range("I113").Value = 4687.5835
$ws.Range("J113").Value = 4475.8335
$ws.Range("K113").Value = 4687.5835
$ws.Range("L113").Value = 4475.8335
$ws.Range("M113").Value = -2517.5835
$ws.Range("N113").Value = -8815.833500000001

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3603.7354
$ws.Range("I132").Value = 2831.65
$ws.Range("J132").Value = 4706.7144
$ws.Range("K132").Value = 8494.950000000001
$ws.Range("L132").Value = 14120.1432
$ws.Range("M132").Value = -5964.950000000001
$ws.Range("N132").Value = -19180.1432

# WVR row 124
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 59476.332
$ws.Range("J124").Value = 59476.332
$ws.Range("L124").Value = 59476.332
$ws.Range("N124").Value = -69296.33199999999

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4169178.5
$ws.Range("I132").Value = 2571.08
$ws.Range("J132").Value = 11113524
$ws.Range("K132").Value = 7713.24
$ws.Range("L132").Value = 33340572
$ws.Range("M132").Value = -5183.24
$ws.Range("N132").Value = -33345632
